$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update column G ("K") values per regenerated save_data (K instead of Strike#)
$kValues = @{
    2 = 0
    4 = 2
    5 = 1
    6 = 0
    7 = 0
    8 = 1
    9 = 2
    10 = 0
    11 = 2
    12 = 1
    13 = 0
    14 = 1
    15 = 0
    16 = 1
    17 = 2
    18 = 0
    19 = 0
    20 = 0
    21 = 1
    22 = 1
    23 = 1
    24 = 0
    25 = 1
    26 = 2
    27 = 0
    28 = 1
    29 = 0
    30 = 0
    31 = 1
    32 = 0
    33 = 0
    34 = 0
    35 = 1
    36 = 0
    37 = 0
    38 = 0
    39 = 1
    40 = 1
    41 = 0
    42 = 1
    43 = 1
    44 = 2
    45 = 0
    46 = 0
    48 = 2
    49 = 0
    50 = 0
    52 = 0
    53 = 1
    54 = 2
    55 = 0
    56 = 0
    57 = 0
    58 = 0
    59 = 2
    60 = 0
    61 = 1
    62 = 1
    63 = 0
    64 = 0
    65 = 2
    66 = 0
    67 = 2
    68 = 1
    69 = 2
    70 = 0
    71 = 2
    72 = 0
    74 = 3
    75 = 1
    76 = 0
    77 = 2
}

foreach ($row in $kValues.Keys) {
    $ws.Cells.Item($row, 7).Value = $kValues[$row]
}
